$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$used = $ws.UsedRange
$rowCount = $used.Rows.Count

for ($r = 1; $r -le $rowCount; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G
    $val = $cell.Value()
    if ($val -ne $null -and $val -ne "") {
        $parts = $val -split ", "
        if ($parts.Count -ge 2) {
            $n = $parts.Count
            $last = $parts[$n - 1]
            $secondLast = $parts[$n - 2]
            $parts[$n - 1] = $secondLast
            $parts[$n - 2] = $last
            $newVal = [string]::Join(", ", $parts)
            if ($newVal -ne $val) {
                $cell.Value = $newVal
            }
        }
    }
}
